$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 updates
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65

# Row 12 updates
$ws.Range("I12").Value = 10
$ws.Range("K12").Value = 2.25
$ws.Range("L12").Value = 9
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 8
$ws.Range("AO12").Value = 6.5
$ws.Range("AW12").Value = 9.5

# Row 13 updates
$ws.Range("I13").Value = 4.85
$ws.Range("K13").Value = 2.12
$ws.Range("L13").Value = 4.9
$ws.Range("X13").Value = 7.4
$ws.Range("AH13").Value = 28
$ws.Range("AP13").Value = 17.5
$ws.Range("AT13").Value = 2.55
$ws.Range("AV13").Value = 65
$ws.Range("AY13").Value = 30
$ws.Range("AZ13").Value = 150
